# Fruta / hortaliza, semanal
# Insert a new weekly record right after row 44 (as new row 45), pushing the
# existing rows 45-73 down to 46-74, then populate the new row with this
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45..73 down to 46..74, leaving a blank row 45 to fill in.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(45, 3).Value = 'Los Lagos'
$ws.Cells.Item(45, 4).Value = 44438
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 'Fruta'
$ws.Cells.Item(45, 7).Value = 100108
$ws.Cells.Item(45, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(45, 9).Value = 100108002
$ws.Cells.Item(45, 10).Value = 'Mango'
$ws.Cells.Item(45, 11).Value = 'Sin especificar'
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 80
$ws.Cells.Item(45, 14).Value = 12000
$ws.Cells.Item(45, 15).Value = 12000
$ws.Cells.Item(45, 16).Value = 12000
$ws.Cells.Item(45, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(45, 18).Value = 'Brasil'
$ws.Cells.Item(45, 19).Value = 3000
$ws.Cells.Item(45, 20).Value = 4
